# Insert a new weekly price record as the second row of observations
# (new row 26) for "Feria Lagunitas de Puerto Montt - Albahaca", shifting
# all subsequent rows down by one (old row 26 -> 27, ..., old row 119 -> 120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 26; this pushes rows 26..119 down to 27..120
# and Excel automatically extends the sheet dimension accordingly.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new observation.
$ws.Range("A26").Value2 = 4
$ws.Range("B26").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C26").Value2 = "Los Lagos"
$ws.Range("D26").Value2 = 44690
$ws.Range("E26").Value2 = 10
$ws.Range("F26").Value2 = 100112052
$ws.Range("G26").Value2 = "Albahaca"
$ws.Range("H26").Value2 = "Sin especificar"
$ws.Range("I26").Value2 = "Primera"
$ws.Range("J26").Value2 = 60
$ws.Range("K26").Value2 = 8000
$ws.Range("L26").Value2 = 8000
$ws.Range("M26").Value2 = 8000
$ws.Range("N26").Value2 = "`$/docena de matas"
$ws.Range("O26").Value2 = "Región Metropolitana"
$ws.Range("P26").Value2 = 1333
$ws.Range("Q26").Value2 = 6
$ws.Range("R26").Value2 = "Hortaliza"
